# Update "想去人数" (F column) figures on the 展览, 演出 and 全部类型 sheets
# to match the newly scraped counts, as published to gh-pages at 456a3b4.

$wb = $excel.ActiveWorkbook

# --- Sheet "展览" ---------------------------------------------------------
$ws1 = $wb.Worksheets.Item("展览")
$sheet1Updates = @{
    2  = 103
    3  = 1250
    4  = 889
    5  = 920
    6  = 1662
    7  = 359
    8  = 1115
    10 = 96
    11 = 241
    12 = 7
    13 = 74
    14 = 595
    15 = 114
    16 = 69
    20 = 61
    21 = 625
    22 = 613
    23 = 103
    26 = 284
    27 = 42
    28 = 15
    29 = 230
}
foreach ($row in $sheet1Updates.Keys) {
    $ws1.Range("F$row").Value = $sheet1Updates[$row]
}

# --- Sheet "演出" ---------------------------------------------------------
$ws2 = $wb.Worksheets.Item("演出")
$sheet2Updates = @{
    7 = 224
    8 = 79
    9 = 609
}
foreach ($row in $sheet2Updates.Keys) {
    $ws2.Range("F$row").Value = $sheet2Updates[$row]
}

# --- Sheet "全部类型" ------------------------------------------------------
$ws4 = $wb.Worksheets.Item("全部类型")
$sheet4Updates = @{
    3  = 103
    4  = 1250
    5  = 889
    6  = 920
    7  = 1662
    8  = 359
    9  = 1115
    12 = 96
    13 = 241
    14 = 7
    15 = 74
    16 = 595
    17 = 114
    18 = 69
    26 = 224
    27 = 224
    28 = 61
    29 = 625
    30 = 613
    31 = 103
    34 = 284
    35 = 79
    36 = 42
    37 = 15
    38 = 230
    39 = 609
}
foreach ($row in $sheet4Updates.Keys) {
    $ws4.Range("F$row").Value = $sheet4Updates[$row]
}
